# #764 adjust cfs note
# Adds the new "SQL Tuesday Sao Paulo 2024 (1099)" attendance row (row 33)
# to the bottom of the tracking table, matching the format of the row
# directly above it (row 32).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the formatting from the last populated row so the new row's date,
# text and percentage cells look consistent with the rest of the table.
$ws.Range("A32:E32").Copy() | Out-Null
$ws.Range("A33").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A33").Value = Get-Date -Year 2024 -Month 12 -Day 17 -Hour 0 -Minute 0 -Second 0
$ws.Range("B33").Value = "SQL Tuesday Sao Paulo 2024 (1099)"
$ws.Range("C33").Value = 152
$ws.Range("D33").Value = 81
$ws.Range("E33").Formula = "=IF(C33=0,0,+(C33-D33)/C33)"

$ws.Range("B33").Select()
